$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date/time number format used for the newly added enquiry rows.
$dateFormat = "yyyy-MM-dd HH:mm:ss"

# Row 13: new enquiry submitted by S1234567A on project 2
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "S1234567A"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "what is going on"
$ws.Range("F13").Value = 45770.229960462966
$ws.Range("F13").NumberFormat = $dateFormat

# Row 14: new enquiry (with reply) submitted by S1234567A on project 3
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "S1234567A"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "test"
$ws.Range("E14").Value = "test"
$ws.Range("F14").Value = 45770.31513549769
$ws.Range("F14").NumberFormat = $dateFormat
$ws.Range("G14").Value = 45770.31800811343
$ws.Range("G14").NumberFormat = $dateFormat
